$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Slides")

# --- Row 13: GRADE S12 ---
# Last_Reviewed becomes a real date/number (was inline text "2026-01-22")
$ws.Range("L13").Value = 46044
$ws.Range("M13").Value = "P1: inserir HR/IC do seguimento 10 anos + figura (paper pendente)"
$ws.Range("N13").Value = "Placeholders adicionados (10y HR/CI + figura/citação). Mantém dados atuais (4.7y) até atualização."

# --- Row 17: GRADE S16 ---
$ws.Range("E17").Value = "READY"
$ws.Range("G17").Value = "Yes"
$ws.Range("H17").Value = "Yes"
$ws.Range("J17").Value = 9
$ws.Range("M17").Value = "P2: revisão final de texto (aula/residentes)"
$ws.Range("N17").Value = "Batch 16–26: layout ok, sem overflow."

# --- Row 18: GRADE S17 ---
$ws.Range("E18").Value = "READY"
$ws.Range("G18").Value = "Yes"
$ws.Range("H18").Value = "Yes"
$ws.Range("J18").Value = 9
$ws.Range("M18").Value = "P2: checar legibilidade da barra MID em projetor"
$ws.Range("N18").Value = "Batch 16–26: padronizado."

# --- Row 19: GRADE S18 ---
$ws.Range("E19").Value = "READY"
$ws.Range("G19").Value = "Yes"
$ws.Range("H19").Value = "Yes"
$ws.Range("J19").Value = 9
$ws.Range("M19").Value = "P2: revisar texto 'rebaxar -1 nível'"
$ws.Range("N19").Value = "Batch 16–26: compactado."

# --- Row 20: GRADE S19 ---
$ws.Range("E20").Value = "READY"
$ws.Range("G20").Value = "Yes"
$ws.Range("H20").Value = "Yes"
$ws.Range("J20").Value = 9
$ws.Range("M20").Value = "P2: validar se 5 domínios RoB2 cabem em telas menores"
$ws.Range("N20").Value = "Compactação de paddings/font para evitar corte do item 5."

# --- Row 21: GRADE S20 ---
$ws.Range("E21").Value = "READY"
$ws.Range("G21").Value = "Yes"
$ws.Range("H21").Value = "Yes"
$ws.Range("J21").Value = 9
$ws.Range("M21").Value = "P2: checar contraste da linha RECOMENDAÇÃO no projetor"
$ws.Range("N21").Value = "Tabela EtD revisada; recomendação destacada."

# --- Row 24: GRADE S24 ---
$ws.Range("E24").Value = "READY"
$ws.Range("G24").Value = "Yes"
$ws.Range("H24").Value = "Yes"
$ws.Range("J24").Value = 9
$ws.Range("M24").Value = "P2: revisar mensagem sobre validação BR (sem datas)"
$ws.Range("N24").Value = "Fonte corrigida para Circulation 2024 PREVENT; removido '2026'."

# --- Row 25: GRADE S25 ---
$ws.Range("E25").Value = "READY"
$ws.Range("G25").Value = "Yes"
$ws.Range("H25").Value = "Yes"
$ws.Range("J25").Value = 9
$ws.Range("M25").Value = "P2: validar narrativa 10y vs 30y"
$ws.Range("N25").Value = "Fontes corrigidas (PREVENT Circulation 2024 + Framingham 30y Circulation 2009)."

# --- Row 26: GRADE S26 ---
$ws.Range("E26").Value = "READY"
$ws.Range("G26").Value = "Yes"
$ws.Range("H26").Value = "Yes"
$ws.Range("J26").Value = 9
$ws.Range("M26").Value = "P2: revisão final da tabela comparativa"
$ws.Range("N26").Value = "Fontes corrigidas (PREVENT, PCE guideline, SCORE2). Mensagem-chave mais cautelosa."

# --- Row 43: GRADE S43 ---
$ws.Range("E43").Value = "READY"
$ws.Range("G43").Value = "Yes"
$ws.Range("H43").Value = "Yes"
$ws.Range("J43").Value = 9
$ws.Range("M43").Value = "P2: checar contraste do painel escuro em projetor"
$ws.Range("N43").Value = "Fix: texto ilegível no navy (Rabdomiólise/contraindicação)."

# --- Row 46: GRADE S46 ---
$ws.Range("E46").Value = "READY"
$ws.Range("G46").Value = "Yes"
$ws.Range("H46").Value = "Yes"
$ws.Range("J46").Value = 9
$ws.Range("M46").Value = "P2: checar contraste do painel BR em projetor"
$ws.Range("N46").Value = "Fix: 'PCSK9i indisponível no SUS' agora em vermelho legível no navy."
